# push - bug fix + feedback update
#
# The backlog-events export ("healthscore_calculado") is refreshed:
#  - the two stale "Evera" (id 148734) rows are dropped
#  - the "Mart Minas" (id 149896) rows move up and get a corrected
#    Healthscore (K) of 10 for the "Metas não atingidas" occurrence
#  - two new FBMDS (id 156244) occurrences are appended
#  - two new Litero (id 156243) occurrences are appended
#  - the trailing Mart Minas (id 154950) row is kept, shifted down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A id Runrunit | B titulo | C estado | D Quadro | E tags
#          F ocorridos | G data | H cliente | I Ranking_de_Eventos
#          J Delta | K Healthscore
$rowsData = @(
    @(149896, "Mart Minas, 15/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Metas não atingidas", "2024-07-15", "Mart Minas", "Metas não atingidas", -2, 10),
    @(149896, "Mart Minas, 15/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Feedback positivo", "2024-07-15", "Mart Minas", "Feedback positivo", 2.5, 10),
    @(156244, "FBMDS, 29/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Planejamento foi reprovado", "2024-07-29", "FBMDS", "Planejamento foi reprovado", -1, 9),
    @(156244, "FBMDS, 29/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Houve atraso nas entregas, isso prejudicou os clientes", "2024-07-29", "FBMDS", "Houve atraso nas entregas, isso prejudicou os clientes", -1, 9),
    @(156243, "Litero, 05/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Cliente pediu proposta", "2024-08-05", "Litero", "Cliente pediu proposta", 2.5, 10),
    @(156243, "Litero, 05/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Feedback positivo", "2024-08-05", "Litero", "Feedback positivo", 2.5, 10),
    @(154950, "Mart Minas, 12/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Houve atraso nas entregas, isso prejudicou os clientes", "2024-08-12", "Mart Minas", "Houve atraso nas entregas, isso prejudicou os clientes", -1, 9)
)

$firstDataRow = 2
$lastRow = $firstDataRow + $rowsData.Length - 1

# Column G ("data") holds plain ISO-looking text (e.g. "2024-07-15") that
# Excel would otherwise auto-coerce into a date serial number. Force the
# whole column range to Text first, write the values, then restore the
# default "General" style so no stray number-format is left behind.
$gRange = $ws.Range("G$firstDataRow`:G$lastRow")
$gRange.NumberFormat = "@"

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $firstDataRow + $i
    $row = $rowsData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
}

$gRange.Style = "Normal"
